$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "52.005.97"
$ws.Range("E2").Value = "  +0.35%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.976.26"
$ws.Range("E3").Value = "  +1.37%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "354.10"
$ws.Range("E5").Value = "  +0.48%  "

# Row 6 - Solana
Set-TextValue "D6" "108.33"
$ws.Range("E6").Value = "  -3.57%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.68%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.84%  "

# Row 10 - Avalanche
Set-TextValue "D10" "38.25"
$ws.Range("E10").Value = "  -2.79%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.47%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.0856"
$ws.Range("E12").Value = "  -4.20%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  -3.51%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "3.445.81"
$ws.Range("E14").Value = "  +1.28%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -2.42%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.985.09"
$ws.Range("E16").Value = "  +1.46%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +1.10%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "52.135.99"
$ws.Range("E18").Value = "  +0.40%  "

# Row 19 - ImmutableX
Set-TextValue "D19" "3.47"
$ws.Range("E19").Value = "  +4.42%  "

# Row 20 - Uniswap
Set-TextValue "D20" "7.49"
$ws.Range("E20").Value = "  -2.06%  "

# Row 21 - InternetComputer (DFINITY)
Set-TextValue "D21" "13.60"
$ws.Range("E21").Value = "  -5.26%  "

# Row 22 - ShibaInu
Set-TextValue "D22" "0.0₃0973"
$ws.Range("E22").Value = "  -1.66%  "

# Row 23 - Litecoin
Set-TextValue "D23" "69.57"
$ws.Range("E23").Value = "  -2.28%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "263.73"
$ws.Range("E24").Value = "  -2.33%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "2.74"
$ws.Range("E25").Value = "  -1.56%  "

# Row 26 - Kaspa
Set-TextValue "D26" "0.178"
$ws.Range("E26").Value = "  +0.02%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -0.31%  "

# Row 28 - Filecoin
Set-TextValue "D28" "7.59"
$ws.Range("E28").Value = "  +2.64%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.108"
$ws.Range("E30").Value = "  -0.49%  "

# Row 31 - Cosmos
Set-TextValue "D31" "10.29"
$ws.Range("E31").Value = "  -2.94%  "

# Row 32 - RenderToken
Set-TextValue "D32" "6.10"
$ws.Range("E32").Value = "  -1.71%  "

# Row 33 - InjectiveProtocol
Set-TextValue "D33" "36.37"
$ws.Range("E33").Value = "  -3.02%  "

# Row 34 - Toncoin
Set-TextValue "D34" "2.20"
$ws.Range("E34").Value = "  -3.20%  "

# Row 35 - OKB
Set-TextValue "D35" "50.80"
$ws.Range("E35").Value = "  -3.90%  "

# Row 36 - VeChain
$ws.Range("E36").Value = "  -2.03%  "

# Row 37 - FirstDigitalUSD
Set-TextValue "D37" "0.999"
$ws.Range("E37").Value = "  +0.02%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "3.20"
$ws.Range("E38").Value = "  -3.38%  "

# Row 39 - Celestia
$ws.Range("E39").Value = "  -4.77%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -3.65%  "

# Row 41 - Stacks
Set-TextValue "D41" "2.71"
$ws.Range("E41").Value = "  +0.58%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  -0.40%  "

# Row 43 - EnergySwap
Set-TextValue "D43" "22.70"
$ws.Range("E43").Value = "  -2.13%  "

# Row 44 - Monero
Set-TextValue "D44" "121.42"
$ws.Range("E44").Value = "  +8.28%  "

# Row 45 - WEMIXToken
$ws.Range("E45").Value = "  -3.21%  "

# Row 46 - Maker
Set-TextValue "D46" "2.119.20"
$ws.Range("E46").Value = "  -2.27%  "

# Row 47 - NEARProtocol
Set-TextValue "D47" "3.36"
$ws.Range("E47").Value = "  -4.73%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").Value = "  -7.38%  "

# Row 49 - TheGraph
$ws.Range("E49").Value = "  -2.49%  "

# Row 50 - BEAM
$ws.Range("E50").Value = "  -2.72%  "

# Row 51 - SEI
Set-TextValue "D51" "0.932"
$ws.Range("E51").Value = "  -0.41%  "
